$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 571.7143
$ws.Range("I28").Value = 496.16666
$ws.Range("J28").Value = 1025
$ws.Range("K28").Value = 496.16666
$ws.Range("L28").Value = 1025
$ws.Range("M28").Value = -11.16665999999998
$ws.Range("N28").Value = -1995

$ws.Range("H33").Value = 536.5
$ws.Range("I33").Value = 245.625
$ws.Range("K33").Value = 245.625
$ws.Range("M33").Value = -16.625

$ws.Range("H41").Value = 1102.5
$ws.Range("I41").Value = 588.3333
$ws.Range("K41").Value = 588.3333
$ws.Range("M41").Value = -148.3333

$ws.Range("H53").Value = 406.5
$ws.Range("I53").Value = 551
$ws.Range("J53").Value = 220.71428
$ws.Range("K53").Value = 551
$ws.Range("L53").Value = 220.71428
$ws.Range("M53").Value = 86
$ws.Range("N53").Value = -1494.71428

$ws.Range("H86").Value = 3265
$ws.Range("J86").Value = 5166.6665
$ws.Range("L86").Value = 5166.6665
$ws.Range("N86").Value = -7412.6665

$ws.Range("H89").Value = 3265
$ws.Range("J89").Value = 5166.6665
$ws.Range("L89").Value = 25833.3325
$ws.Range("N89").Value = -37065.3325

$ws.Range("H98").Value = 3845.0715
$ws.Range("I98").Value = 3669.5
$ws.Range("K98").Value = 3669.5
$ws.Range("M98").Value = -2171.5

$ws.Range("H111").Value = 1492.7142
$ws.Range("I111").Value = 1208.1666
$ws.Range("K111").Value = 3624.4998
$ws.Range("M111").Value = -557.4998000000001

$ws.Range("H113").Value = 6279.8
$ws.Range("I113").Value = 3799.6667
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 3799.6667
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -545.6667000000002
$ws.Range("N113").Value = -16508

$ws.Range("H116").Value = 5435.9287
$ws.Range("J116").Value = 7534.75
$ws.Range("L116").Value = 7534.75
$ws.Range("N116").Value = -14418.75

$ws.Range("H122").Value = 3845.0715
$ws.Range("I122").Value = 3669.5
$ws.Range("K122").Value = 11008.5
$ws.Range("M122").Value = -8558.5

$ws.Range("H129").Value = 3502.353
$ws.Range("I129").Value = 1015.2
$ws.Range("K129").Value = 3045.6
$ws.Range("M129").Value = 1954.4

$ws.Range("H131").Value = 2269.3
$ws.Range("J131").Value = 2500
$ws.Range("L131").Value = 7500
$ws.Range("N131").Value = -17580

$ws.Range("H137").Value = 2993.5557
$ws.Range("J137").Value = 3770.5
$ws.Range("L137").Value = 11311.5
$ws.Range("N137").Value = -16411.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1594

$ws.Range("H37").Value = 23750
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19727

$ws.Range("H55").Value = 18333.334
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 25000
$ws.Range("K55").Value = 5000
$ws.Range("L55").Value = 25000
$ws.Range("M55").Value = -4685
$ws.Range("N55").Value = -25630

$ws.Range("H132").Value = 3788.9583
$ws.Range("I132").Value = 3632.625
$ws.Range("J132").Value = 4101.625
$ws.Range("K132").Value = 10897.875
$ws.Range("L132").Value = 12304.875
$ws.Range("M132").Value = -8367.875
$ws.Range("N132").Value = -17364.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1030.2
$ws.Range("I64").Value = 1476.3334
$ws.Range("J64").Value = 732.7778
$ws.Range("K64").Value = 1476.3334
$ws.Range("L64").Value = 732.7778
$ws.Range("M64").Value = -1251.3334
$ws.Range("N64").Value = -1182.7778

$ws.Range("H67").Value = 1030.2
$ws.Range("I67").Value = 1476.3334
$ws.Range("J67").Value = 732.7778
$ws.Range("K67").Value = 1476.3334
$ws.Range("L67").Value = 732.7778
$ws.Range("M67").Value = -696.3334
$ws.Range("N67").Value = -2292.7778

$ws.Range("H94").Value = 1847.591
$ws.Range("J94").Value = 1605
$ws.Range("L94").Value = 1605
$ws.Range("N94").Value = -2507

$ws.Range("H105").Value = 2752.7058
$ws.Range("I105").Value = 2316.2222
$ws.Range("K105").Value = 2316.2222
$ws.Range("M105").Value = -569.2222000000002

$ws.Range("H107").Value = 1857
$ws.Range("J107").Value = 3750
$ws.Range("L107").Value = 3750
$ws.Range("N107").Value = -7590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3424.5
$ws.Range("I31").Value = 2596.0908
$ws.Range("J31").Value = 5247
$ws.Range("K31").Value = 2596.0908
$ws.Range("L31").Value = 5247
$ws.Range("M31").Value = -2301.0908
$ws.Range("N31").Value = -5837

$ws.Range("H34").Value = 3424.5
$ws.Range("I34").Value = 2596.0908
$ws.Range("J34").Value = 5247
$ws.Range("K34").Value = 2596.0908
$ws.Range("L34").Value = 5247
$ws.Range("M34").Value = -2394.0908
$ws.Range("N34").Value = -5651

$ws.Range("H94").Value = 2948.75
$ws.Range("J94").Value = 2948
$ws.Range("L94").Value = 2948
$ws.Range("N94").Value = -3850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 500.89474
$ws.Range("I5").Value = 498.2
$ws.Range("J5").Value = 511
$ws.Range("K5").Value = 1494.6
$ws.Range("L5").Value = 1533
$ws.Range("M5").Value = -1382.6
$ws.Range("N5").Value = -1757

$ws.Range("H80").Value = 2500.25
$ws.Range("I80").Value = 1999.5
$ws.Range("J80").Value = 3001
$ws.Range("K80").Value = 5998.5
$ws.Range("L80").Value = 9003
$ws.Range("M80").Value = -5062.5
$ws.Range("N80").Value = -10875

$ws.Range("H83").Value = 2500.25
$ws.Range("I83").Value = 1999.5
$ws.Range("J83").Value = 3001
$ws.Range("K83").Value = 17995.5
$ws.Range("L83").Value = 27009
$ws.Range("M83").Value = -13315.5
$ws.Range("N83").Value = -36369

$ws.Range("H120").Value = 21190
$ws.Range("J120").Value = 28571.428
$ws.Range("L120").Value = 85714.284
$ws.Range("N120").Value = -95390.284

$ws.Range("H122").Value = 1310.1875
$ws.Range("I122").Value = 1191.375
$ws.Range("J122").Value = 1429
$ws.Range("K122").Value = 10722.375
$ws.Range("L122").Value = 12861
$ws.Range("M122").Value = -8272.375
$ws.Range("N122").Value = -17761

$ws.Range("H132").Value = 1762.8
$ws.Range("I132").Value = 1025.8
$ws.Range("J132").Value = 2499.8
$ws.Range("K132").Value = 9232.199999999999
$ws.Range("L132").Value = 22498.2
$ws.Range("M132").Value = -6702.199999999999
$ws.Range("N132").Value = -27558.2

$ws.Range("H135").Value = 500.89474
$ws.Range("I135").Value = 498.2
$ws.Range("J135").Value = 511
$ws.Range("K135").Value = 4483.8
$ws.Range("L135").Value = 4599
$ws.Range("M135").Value = -1948.8
$ws.Range("N135").Value = -9669

$ws.Range("H137").Value = 6843
$ws.Range("J137").Value = 7549.143
$ws.Range("L137").Value = 22647.429
$ws.Range("N137").Value = -32847.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5401.2
$ws.Range("I122").Value = 3402.4
$ws.Range("J122").Value = 7400
$ws.Range("K122").Value = 10207.2
$ws.Range("L122").Value = 22200
$ws.Range("M122").Value = -7757.200000000001
$ws.Range("N122").Value = -27100

$ws.Range("H136").Value = 3436
$ws.Range("I136").Value = 3106.7693
$ws.Range("J136").Value = 3911.5557
$ws.Range("K136").Value = 9320.3079
$ws.Range("L136").Value = 11734.6671
$ws.Range("M136").Value = -6770.3079
$ws.Range("N136").Value = -16834.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50000
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50462

$ws.Range("H126").Value = 2099.75
$ws.Range("I126").Value = 2099.75
$ws.Range("K126").Value = 6299.25
$ws.Range("M126").Value = -3829.25

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070

Write-Host "Applied profit updates"
